# Applies the update described by the diff: a corrected Taxonsorteringsordning
# (B2), and the species records in rows 3 & 4 swapping places (with revised
# B-values) plus row 3 gaining a Biotop/Biotop-beskrivning pair that row 4
# loses.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: only the Taxonsorteringsordning changes ---
$ws.Range("B2").Value = 76733

# --- Row 3: becomes the "Ögonpyrola" record (previously in row 4) ---
$ws.Range("A3").Value = 112389108
$ws.Range("B3").Value = 103766
$ws.Range("E3").Value = 221725
$ws.Range("F3").Value = "Ögonpyrola"
$ws.Range("G3").Value = "Moneses uniflora"
$ws.Range("H3").Value = "(L.) A. Gray"

# Antal / dates are stored as text in this sheet, so force the Text format
# before writing to stop Excel from reinterpreting "50" as a number or the
# dates as date serials.
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "50"
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2023-07-12"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2023-07-12"

# Row 3 gains Biotop / Biotop-beskrivning values
$ws.Range("AH3").Value = "Frisk gräsmark med lövträd"
$ws.Range("AI3").Value = "slåtterängsmark"

# --- Row 4: becomes the "Skogsnattviol" record (previously in row 3) ---
$ws.Range("A4").Value = 112389257
$ws.Range("B4").Value = 96770
$ws.Range("E4").Value = 223621
$ws.Range("F4").Value = "Skogsnattviol"
$ws.Range("G4").Value = "Platanthera bifolia subsp. latiflora"
$ws.Range("H4").Value = "(Drejer) Løjtnant"

$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "24"
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2023-06-18"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2023-06-18"

# Row 4 loses its Biotop / Biotop-beskrivning values entirely
$ws.Range("AH4:AI4").ClearContents()
